$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados..." timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 02:10"

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1768116
$ws.Range("C4").Value = 22313
$ws.Range("D4").Value = 498681
$ws.Range("E4").Value = 1166115
$ws.Range("G4").Value = 1213
$ws.Range("H4").Value = 103320

# Row 5 (Brasil)
$ws.Range("B5").Value = 438812
$ws.Range("C5").Value = 24151
$ws.Range("D5").Value = 193181
$ws.Range("E5").Value = 218640
$ws.Range("G5").Value = 1294
$ws.Range("H5").Value = 26991

# Row 49 (Panama)
$ws.Range("B49").Value = 12131
$ws.Range("C49").Value = 403
$ws.Range("E49").Value = 4432
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 320

# Row 124 (Uruguay)
$ws.Range("B124").Value = 811
$ws.Range("C124").Value = 8
$ws.Range("D124").Value = 654
$ws.Range("E124").Value = 135

# Row 165 (Brunei)
$ws.Range("B165").Value = 149
$ws.Range("C165").Value = 17
$ws.Range("D165").Value = 28
$ws.Range("E165").Value = 117
$ws.Range("H165").Value = 4

# Row 166 (Islas Caimanes)
$ws.Range("B166").Value = 141
$ws.Range("D166").Value = 138
$ws.Range("E166").Value = 1
$ws.Range("H166").Value = 2

# Row 167 (Bermudas)
$ws.Range("B167").Value = 140
$ws.Range("D167").Value = 67
$ws.Range("E167").Value = 72
$ws.Range("H167").Value = 1

# Row 168 (Zimbabue)
$ws.Range("B168").Value = 140
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 92
$ws.Range("E168").Value = 39
$ws.Range("H168").Value = 9

# Row 172 (Libia)
$ws.Range("D172").Value = 41
$ws.Range("E172").Value = 59

# Row 173 (Aruba)
$ws.Range("D173").Value = 98
$ws.Range("E173").Value = 0

# Row 210 (Montserrat)
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# Row 211 (Seychelles)
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
